$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows covering 2021-05-28 through 2021-06-28 (serials 44344-44375)
$newRows = @(
    @(270, 44344, 0, 3, 80.29978586723769),
    @(271, 44345, 0, 2, 53.53319057815846),
    @(272, 44346, 0, 1, 26.76659528907923),
    @(273, 44347, 0, 0, 0),
    @(274, 44348, 0, 0, 0),
    @(275, 44349, 0, 0, 0),
    @(276, 44350, 0, 0, 0),
    @(277, 44351, 0, 0, 0),
    @(278, 44352, 0, 0, 0),
    @(279, 44353, 0, 0, 0),
    @(280, 44354, 0, 0, 0),
    @(281, 44355, 0, 0, 0),
    @(282, 44356, 0, 0, 0),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 0, 0, 0),
    @(285, 44359, 0, 0, 0),
    @(286, 44360, 0, 0, 0),
    @(287, 44361, 0, 0, 0),
    @(288, 44362, 0, 0, 0),
    @(289, 44363, 0, 0, 0),
    @(290, 44364, 0, 0, 0),
    @(291, 44365, 0, 0, 0),
    @(292, 44366, 0, 0, 0),
    @(293, 44367, 0, 0, 0),
    @(294, 44368, 0, 0, 0),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 0, 0, 0),
    @(297, 44371, 0, 0, 0),
    @(298, 44372, 0, 0, 0),
    @(299, 44373, 0, 0, 0),
    @(300, 44374, 0, 0, 0),
    @(301, 44375, 0, 0, 0)
)

$lastExistingRow = 269
$firstNewRow = 270
$lastNewRow = 301

# Copy the formatting (style) of the last existing data row (A269) down across
# the new rows in column A, so the date style (s="2") is applied consistently.
$ws.Range("A" + $lastExistingRow).Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0

